# hotfix double-encoding-issue by using triple-braces
#
# The merge-field placeholders in this offer-letter template were written
# with double curly braces ("{{Field__c}}"). The templating engine that
# consumes this workbook double-encodes those, so every placeholder is
# rewritten here with triple braces ("{{{Field__c}}}"). Only the literal
# text of the placeholder cells changes - labels / static Japanese text
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Opening paragraph: "{{AccountName__c}}(..." -> "{{{AccountName__c}}}(..."
$ws.Range("B4").Value = "{{{AccountName__c}}}（以下、甲という。）と　株式会社サンプル（以下、乙という。）は、"

# 氏名 (employer name) merge field
$ws.Range("C7").Value = "{{{AccountName__c}}}"

# 現住所 (employer address) merge field
$ws.Range("C8").Value = "{{{AccountAddress__c}}}"

# 雇用期間 (employment period)
$ws.Range("C11").Value = "{{{StartDateFormat__c}}} 〜 {{{EndDateFormat__c}}}  "

# 勤務場所 (work location)
$ws.Range("C12").Value = "{{{Address__c}}}"

# 仕事内容 (job description)
$ws.Range("C13").Value = "{{{JobDescription__c}}} "

# 就業時間 (working hours)
$ws.Range("C14").Value = "{{{StartTime__c}}} 〜 {{{EndTime__c}}}  "

# 所定外労働の有無 (overtime)
$ws.Range("C15").Value = "{{{hasOverTime__c}}} "

# 休暇 (holidays)
$ws.Range("C16").Value = "{{{HoliDayType__c}}} "

# 賃金 (salary)
$ws.Range("C17").Value = "基本給(月)　{{{Salary__c}}}万円"

# 賃金締切日 (pay period cutoff date)
$ws.Range("C18").Value = "{{{DueDate__c}}} "

# 賃金支払日 (pay date)
$ws.Range("C19").Value = "{{{SalaryDate__c}}} "

# Drop the stray "H3" cell selection that had been saved with the sheet
# (the template should open with the default/home selection instead).
[void]$ws.Range("A1").Select()

